{"js": "// Add two new paragraphs at the end of the document body, after the\n// existing \"My name is Vignesh\" paragraph:\n//   \"From Manjeshwar\"\n//   \"Studying in Vivekananda Polytechnic Puttur\"\n\nconst body = context.document.body;\nbody.paragraphs.load(\"text\");\nawait context.sync();\n\n// Insert the first new paragraph at the very end of the body (after the\n// last existing paragraph, \"My name is Vignesh\").\nconst p1 = body.insertParagraph(\"From Manjeshwar\", Word.InsertLocation.end);\nawait context.sync();\n\n// Insert the second new paragraph right after the one we just added.\np1.insertParagraph(\"Studying in Vivekananda Polytechnic Puttur\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Add two new paragraphs at the end of the document body, after the\n# existing \"My name is Vignesh\" paragraph:\n#   \"From Manjeshwar\"\n#   \"Studying in Vivekananda Polytechnic Puttur\"\n\n$d = $word.ActiveDocument\n\n# Insert a new (empty) paragraph right after the current last paragraph,\n# then fill it in with the first new line of text.\n$last = $d.Paragraphs.Last\n$last.Range.InsertParagraphAfter()\n$p1 = $d.Paragraphs.Last\n$p1.Range.Text = \"From Manjeshwar\"\n\n# Repeat for the second new line, appended after the one we just added.\n$p1 = $d.Paragraphs.Last\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Last\n$p2.Range.Text = \"Studying in Vivekananda Polytechnic Puttur\"\n"}
